$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.571.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.76%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.643.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +4.42%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.9981"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.51%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'307.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.92%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.9987"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.30%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.3786"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.29%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'53.08"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +6.45%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.3697"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +4.37%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +6.13%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.08204"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +3.19%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.9992"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.43%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'23.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +7.33%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.680"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +4.30%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.00001291"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +5.88%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'7.487"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.03%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.640.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +3.98%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'95.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +3.78%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.06960"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.44%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'18.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +5.02%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.609"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +4.02%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.9983"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.34%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'23.565.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.77%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'13.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +3.74%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'3.134"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +12.08%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.423"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.23%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'21.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +4.73%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'151.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +3.22%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'5.336"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +3.64%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'136.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +4.00%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'2.428"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +4.36%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'6.869"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +5.60%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.816.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +3.73%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.9794"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +5.43%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.02834"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +7.90%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'10.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +5.54%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.07507"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.36%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'6.248"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +4.86%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.2549"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.73%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +1.44%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +4.61%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.7199"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'12.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +7.63%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +10.30%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.6665"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +5.64%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.373"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'4.044"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.01%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.9978"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.29%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.08077"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.03%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'131.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.50%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.217"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +3.19%  "
$ws.Range("E51").Style = "Normal"
